# Update "Förändrad" (column C) date serials from 45190 to 45192 for all data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C246").Value = 45192
